$wb = $excel.ActiveWorkbook

# --- 1. Update the "Date" value on the Metadata sheet ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value2 = "2024-03-19T13:17:15+00:00"

# --- 2. Swap the two "Mapping" columns (header + data + column widths) on Elements sheet ---
$elements = $wb.Worksheets.Item("Elements")

# Swap header labels AK1 / AL1 (RIM Mapping <-> Spécification métier mapping)
$akHeader = $elements.Range("AK1").Value2
$alHeader = $elements.Range("AL1").Value2
$elements.Range("AK1").Value2 = $alHeader
$elements.Range("AL1").Value2 = $akHeader

# Swap the data values for rows 2-6 in columns AK / AL
for ($r = 2; $r -le 6; $r++) {
    $akCell = $elements.Cells.Item($r, 37)
    $alCell = $elements.Cells.Item($r, 38)
    $akVal = $akCell.Value2
    $alVal = $alCell.Value2
    $akCell.Value2 = $alVal
    $alCell.Value2 = $akVal
}

# Swap the column widths of AK (37) and AL (38)
$elements.Columns.Item(37).ColumnWidth = 65.17578125
$elements.Columns.Item(38).ColumnWidth = 24.98046875
